# Update "想去人数" (interest count) values in column F across the four
# worksheets of the workbook, matching the data refresh captured in the
# commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 232
$ws1.Range("F4").Value  = 819
$ws1.Range("F5").Value  = 238
$ws1.Range("F6").Value  = 402
$ws1.Range("F7").Value  = 564
$ws1.Range("F9").Value  = 60
$ws1.Range("F10").Value = 336
$ws1.Range("F11").Value = 130
$ws1.Range("F12").Value = 611
$ws1.Range("F13").Value = 80
$ws1.Range("F14").Value = 1764
$ws1.Range("F15").Value = 324
$ws1.Range("F16").Value = 2457
$ws1.Range("F17").Value = 299

# Sheet 2: 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F13").Value = 85

# Sheet 3: 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5286
$ws3.Range("F3").Value = 307
$ws3.Range("F4").Value = 200

# Sheet 4: 全部类型 (All types) - aggregate sheet, mirrors the three above
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 5286
$ws4.Range("F4").Value  = 307
$ws4.Range("F6").Value  = 200
$ws4.Range("F7").Value  = 232
$ws4.Range("F13").Value = 819
$ws4.Range("F16").Value = 238
$ws4.Range("F17").Value = 402
$ws4.Range("F18").Value = 564
$ws4.Range("F20").Value = 60
$ws4.Range("F22").Value = 336
$ws4.Range("F23").Value = 130
$ws4.Range("F26").Value = 611
$ws4.Range("F27").Value = 80
$ws4.Range("F28").Value = 85
$ws4.Range("F29").Value = 1764
$ws4.Range("F30").Value = 324
$ws4.Range("F31").Value = 2458
$ws4.Range("F33").Value = 299
